$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some price values round-trip through Excel's numeric auto-detection and
# lose a significant trailing zero (e.g. "1.00" -> 1, "23.40" -> 23.4).
# Force those specific cells to be stored as text so the displayed value
# keeps its original trailing zeros, matching the source data feed.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.976.17"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.453.60"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "580.84"
$ws.Range("E5").Value = "  +0.35%  "

# Row 6 - Solana
$ws.Range("D6").Value = "150.95"
$ws.Range("E6").Value = "  +2.60%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.488"
$ws.Range("E8").Value = "  +1.52%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "8.08"
$ws.Range("E9").Value = "  +6.34%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +4.23%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.040.23"
$ws.Range("E12").Value = "  -0.60%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.36%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "28.48"
$ws.Range("E14").Value = "  -4.25%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.459.11"
$ws.Range("E15").Value = "  -0.29%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.46%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.960.52"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +2.22%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "14.63"
$ws.Range("E19").Value = "  +1.92%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "9.05"
$ws.Range("E20").Value = "  -1.79%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "389.32"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.97%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "75.21"
$ws.Range("E23").Value = "  +0.87%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.591.50"
$ws.Range("E25").Value = "  -0.54%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -0.64%  "

# Row 27 - Kaspa
$ws.Range("D27").Value = "0.187"
$ws.Range("E27").Value = "  +4.44%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "7.76"
$ws.Range("E28").Value = "  +2.60%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.12%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "8.06"
$ws.Range("E30").Value = "  -1.15%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.28%  "

# Row 32 - USDe
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33 - Fetch.AI
$ws.Range("D33").Value = "1.36"
$ws.Range("E33").Value = "  -1.95%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "23.40"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").Value = "  +3.59%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +4.66%  "

# Row 37 - EnergySwap
$ws.Range("D37").Value = "31.86"
$ws.Range("E37").Value = "  +0.88%  "

# Row 38 - Aptos
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  -1.43%  "

# Row 39 - Monero
$ws.Range("D39").Value = "169.21"
$ws.Range("E39").Value = "  +0.01%  "

# Row 40 - RenzoRestakedETH
$ws.Range("D40").Value = "3.489.20"
$ws.Range("E40").Value = "  -0.64%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0787"
$ws.Range("E41").Value = "  +3.17%  "

# Row 42 - was OKB, now Mantle
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.787"
$ws.Range("E42").Value = "  -1.45%  "

# Row 43 - was Mantle, now OKB
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "42.80"
$ws.Range("E43").Value = "  +1.24%  "

# Row 44 - Stacks
$ws.Range("D44").Value = "1.72"
$ws.Range("E44").Value = "  -0.04%  "

# Row 45 - Filecoin
$ws.Range("D45").Value = "4.42"
$ws.Range("E45").Value = "  -1.37%  "

# Row 46 - ONDO
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -0.98%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.559.14"
$ws.Range("E47").Value = "  -1.82%  "

# Row 48 - Cosmos
$ws.Range("D48").Value = "6.97"
$ws.Range("E48").Value = "  +3.23%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "2.25"
$ws.Range("E49").Value = "  +0.08%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "23.01"
$ws.Range("E50").Value = "  +0.04%  "

# Row 51 - FirstDigitalUSD
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.34%  "
